$d = $word.ActiveDocument

# Helper: replace the full text of a paragraph while preserving any leading
# empty run (e.g. <w:r/>) that precedes the text-bearing run, by using
# InsertXML on a range that spans only the existing text run(s) (an empty
# run contributes zero characters to Range.Text, so it falls outside the
# replaced span).
function Set-ParagraphText {
    param($doc, [string]$oldText, [string]$newText)
    foreach ($p in $doc.Paragraphs) {
        $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $oldText) {
            $r = $doc.Range($p.Range.Start, $p.Range.End - 1)
            $escaped = $newText -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'
            $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p><w:r><w:t>' + $escaped + '</w:t></w:r></w:p></w:body></w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'
            $r.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# 1. Title (Heading1 at top, and the bold run near the end share the same text)
$d.Content.Find.Execute(
    "Play BountyPop Free: Exciting PopWins and Bonus Features", $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Play BountyPop Free - Review of the Exciting Online Slot Game", 2)

# 2. "What we like" bullet list rewrites
Set-ParagraphText $d "PopWins feature creates exciting gameplay" "PopWins feature"
Set-ParagraphText $d "High volatility and 96% RTP can lead to big winnings" "Game volatility and RTP"
Set-ParagraphText $d "Multiplier Wheel and Wheel of Bets add extra chances for bonuses" "Multiplier wheel"
Set-ParagraphText $d "Characters inspire mystery, adventure, and magic" "Wheel of Bets"

# 3. "What we don't like" bullet list: rewrite first item, remove second entirely
Set-ParagraphText $d "Graphics and sound design are lackluster" "Graphics and sound design"

foreach ($p in $d.Paragraphs) {
    $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "Wheel of Bets can lead to disappointment if losing") {
        $p.Range.Delete()
        break
    }
}

# 4. Meta description (italic) text rewrite
$d.Content.Find.Execute(
    "Read our review of BountyPop and play for free. Experience exciting PopWins feature, high volatility, and bonus games like Multiplier Wheel and Wheel of Bets.",
    $true, $true, $false, $false, $false,
    $true, 1, $false,
    "Read our review of BountyPop to discover the thrilling features and play for free.", 2)
